# 自动更新价格数据: insert a new top data row (row 2) for 2026-01-15,
# pushing every existing data row down by one and carrying the last
# row's values into the newly vacated bottom row (matches the source
# site's "repeat last price" pattern seen across the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room at the bottom of the table (row 57) first. Row 56 has no
# special formatting, so the freshly-inserted row inherits none either
# (inserting directly above row 2 would instead have pulled down the
# bold header formatting from row 1).
$ws.Rows.Item(57).Insert()

# Shift every existing data row down by one, bottom row first so we
# never overwrite a row before it has been copied.
for ($r = 56; $r -ge 2; $r--) {
    $dateText = $ws.Cells.Item($r, 1).Text
    $b = $ws.Cells.Item($r, 2).Text
    $c = $ws.Cells.Item($r, 3).Text
    $d = $ws.Cells.Item($r, 4).Text

    # Date column holds plain text, not a real date value - force text
    # formatting before writing so Excel doesn't reinterpret the string
    # as a date serial, then clear the formatting override so the cell
    # ends up styled exactly like its neighbours (no explicit style).
    $ws.Cells.Item($r + 1, 1).NumberFormat = "@"
    $ws.Cells.Item($r + 1, 1).Value = $dateText
    $ws.Cells.Item($r + 1, 1).ClearFormats()

    $ws.Cells.Item($r + 1, 2).Value = $b
    $ws.Cells.Item($r + 1, 3).Value = $c
    $ws.Cells.Item($r + 1, 4).Value = $d
}

# New top row: next day's date with the same commodity values as the
# rest of the series.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-01-15"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
